$d = $word.ActiveDocument

# The document's footers/header carry two Pearson logo pictures (currently
# named "image1.png") and one BTec logo picture (currently named
# "image2.jpg"). This commit swaps those display names:
#   Pearson logo pictures : image1.png -> image2.png
#   BTec logo picture     : image2.jpg -> image1.jpg
#
# InlineShape has no settable Name property in the Word object model, so
# each picture is temporarily converted to a floating Shape (which does
# expose Name), renamed, then converted back to an inline picture.

function Rename-InlineLogo($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShp = $range.InlineShapes.Item($i)
        $floatingShp = $inlineShp.ConvertToShape()
        $floatingShp.Name = $newName
        $floatingShp.ConvertToInlineShape()
    }
}

foreach ($sec in $d.Sections) {
    for ($hfIdx = 1; $hfIdx -le 3; $hfIdx++) {
        $ftr = $sec.Footers.Item($hfIdx)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            Rename-InlineLogo $ftr.Range "image2.png"
        }

        $hdr = $sec.Headers.Item($hfIdx)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            Rename-InlineLogo $hdr.Range "image1.jpg"
        }
    }
}
